$d = $word.ActiveDocument
$xml = @"
<w:p>
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
      </w:pPr>
      <w:r>
        <w:t>Description of the Problem</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
      </w:pPr>
      <w:r>
        <w:t>Outline of problem.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>I intend to design a program to find paths between webpages. The program will include the following features: a help screen, text based ui to help user find their way, a database of webpages crawled</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> that will be created using SQL</w:t>
      </w:r>
      <w:r>
        <w:t>, with a linked database of pages they have linked to</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> and a pathfinding system</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>The end users of my program will be people willing to find orphan links in webpages</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> if a search system is not used</w:t>
      </w:r>
      <w:r>
        <w:t>, therefore will probably be tech literate.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">My project meets the advanced higher </w:t>
      </w:r>
      <w:r>
        <w:t>computing requirements as it will have a UI suitable for tech literate users with validation for if the pages have valid urls by using a try catch with a get() procedure and checking the code sent from the sever</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> is not 404 and that the domain exists. My project will interface with an SQL database, creating a database and </w:t>
      </w:r>
      <w:r>
        <w:t>writing and reading URL’s from it.</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">A piece of software that will </w:t>
      </w:r>
      <w:r>
        <w:t>b</w:t>
      </w:r>
      <w:r>
        <w:t>e</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> able to crawl webpages. </w:t>
      </w:r>
      <w:r>
        <w:t>The webpages crawled and the pages they link to should be stored in a database</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> including the link to them using SQL</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>The database can then be read and a path can be found.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> And then the path is displayed to the user.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Scope the clearly defined outline of what the solution will deliver in terms of functionality</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Boundaries: are the limitations of the project</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Constraints the restrictions that apply to the development.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:sectPr>
          <w:pgSz w:w="11906" w:h="16838"/>
          <w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="708" w:footer="708" w:gutter="0"/>
          <w:cols w:space="708"/>
          <w:docGrid w:linePitch="360"/>
        </w:sectPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
      </w:pPr>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>UML</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:sectPr>
          <w:pgSz w:w="16838" w:h="11906" w:orient="landscape"/>
          <w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="708" w:footer="708" w:gutter="0"/>
          <w:cols w:space="708"/>
          <w:docGrid w:linePitch="360"/>
        </w:sectPr>
      </w:pPr>
    </w:p>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
      </w:pPr>
      <w:r>
        <w:t>Requirements</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
      </w:pPr>
      <w:r>
        <w:t>End user</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> requirements</w:t>
      </w:r>
      <w:r>
        <w:t>:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">User must be able to use </w:t>
      </w:r>
      <w:r>
        <w:t>text-based</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> UI to input a starting website, and an end website to find a path to.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">User must be able to input </w:t>
      </w:r>
      <w:r>
        <w:t>the</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> number of moves they wish it to be done in.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>The user must be able to view the requested path, or receive an error message that there is no path</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
      </w:pPr>
      <w:r>
        <w:t>Functional requirements:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>T</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">he program must display a UI </w:t>
      </w:r>
      <w:r>
        <w:t>that can take in a starting page and end page with number of moves.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">The program will </w:t>
      </w:r>
      <w:r>
        <w:t>be able to crawl a url and find all links on the url, follow them, and repeat the process</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> until the maximum jumps is achieved</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">The webpage’s url is to be stored then </w:t>
      </w:r>
      <w:r>
        <w:t>all</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> the links leading off also need to be stored</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> in a database using</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> SQL</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>These are to be stored in a database using sql</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> queries.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>The program will then call the separate path finding algorithm</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>SQL</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> queries then need to be written into a 2d array.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>(possible: sort 2d array so that easier to read.</w:t>
      </w:r>
      <w:r>
        <w:t>)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>Then a node map object is created using the data in the 2d array.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Then using a pathfinding algorithm, </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">find </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>A</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> path (not shortest) from first link to second link</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> in the maximum number of jumps given</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading1"/>
      </w:pPr>
      <w:r>
        <w:t>Project plan</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Heading2"/>
      </w:pPr>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:t>Identified tasks:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
      </w:r>
      <w:r>
        <w:tab/>
      </w:r>
    </w:p>
    
"@
$d.Content.InsertXML($xml)
